{"js": "// The tutorial-assignment bullet list inside \"Workshops: multimedia tutorial\"\n// had a closing bullet point (\"Wrap up the workshop by re-stating the main\n// principles and highlighting examples of the work of the participants.\")\n// that duplicated/contradicted the rest of the assignment instructions.\n// Remove that entire bulleted paragraph, leaving the preceding \"...of time.\"\n// bullet followed directly by the \"Please post your tutorial...\" paragraph.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"Wrap up the workshop by re-stating the main principles\";\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(marker) !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the 'Wrap up the workshop...' paragraph to remove.\");\n}\n\ntarget.delete();\nawait context.sync();\n", "ps1": "# The tutorial-assignment bullet list inside \"Workshops: multimedia tutorial\"\n# ends with a bullet point that duplicated/contradicted the rest of the\n# assignment instructions:\n#   \"Wrap up the workshop by re-stating the main principles and highlighting\n#    examples of the work of the participants.\"\n# Remove that entire bulleted paragraph so the \"...of time.\" bullet is\n# followed directly by the \"Please post your tutorial...\" paragraph.\n\n$d = $word.ActiveDocument\n\n$marker = \"Wrap up the workshop by re-stating the main principles\"\n$target = $null\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*$marker*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"Could not locate the 'Wrap up the workshop...' paragraph to remove.\"\n}\n\n# Deleting the paragraph's Range (which includes its paragraph mark) removes\n# the whole bulleted paragraph, so the following paragraph shifts up cleanly.\n$target.Range.Delete()\n"}
